$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1138.5
$ws.Range("I21").Value = 1138.5
$ws.Range("K21").Value = 1138.5
$ws.Range("M21").Value = -670.5
$ws.Range("H23").Value = 1138.5
$ws.Range("I23").Value = 1138.5
$ws.Range("K23").Value = 1138.5
$ws.Range("M23").Value = -904.5
$ws.Range("H40").Value = 2600
$ws.Range("I40").Value = 2250
$ws.Range("K40").Value = 2250
$ws.Range("M40").Value = -2075
$ws.Range("H43").Value = 3739.8
$ws.Range("I43").Value = 3399
$ws.Range("K43").Value = 3399
$ws.Range("M43").Value = -3330
$ws.Range("H70").Value = 148934.4
$ws.Range("I70").Value = 1850
$ws.Range("J70").Value = 246990.67
$ws.Range("K70").Value = 5550
$ws.Range("L70").Value = 740972.01
$ws.Range("M70").Value = -5280
$ws.Range("N70").Value = -741512.01
$ws.Range("H73").Value = 148934.4
$ws.Range("I73").Value = 1850
$ws.Range("J73").Value = 246990.67
$ws.Range("K73").Value = 5550
$ws.Range("L73").Value = 740972.01
$ws.Range("M73").Value = -4614
$ws.Range("N73").Value = -742844.01
$ws.Range("H103").Value = 1900
$ws.Range("I103").Value = 1900
$ws.Range("K103").Value = 5700
$ws.Range("M103").Value = -5114
$ws.Range("H106").Value = 32940.35
$ws.Range("I106").Value = 34461.23
$ws.Range("K106").Value = 34461.23
$ws.Range("M106").Value = -33830.23
$ws.Range("H113").Value = 4562.75
$ws.Range("I113").Value = 4126.25
$ws.Range("K113").Value = 4126.25
$ws.Range("M113").Value = -872.25
$ws.Range("H138").Value = 4264.2793
$ws.Range("I138").Value = 3443.4546
$ws.Range("J138").Value = 5038.2
$ws.Range("K138").Value = 10330.3638
$ws.Range("L138").Value = 15114.6
$ws.Range("M138").Value = -5190.363799999999
$ws.Range("N138").Value = -25394.6
$ws.Range("H141").Value = 1501.7307
$ws.Range("I141").Value = 1501.7307
$ws.Range("K141").Value = 4505.1921
$ws.Range("M141").Value = 674.8078999999998
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 34996.668
$ws.Range("J24").Value = 34996.668
$ws.Range("L24").Value = 34996.668
$ws.Range("N24").Value = -35744.668
$ws.Range("H32").Value = 5822.4863
$ws.Range("I32").Value = 3498.121
$ws.Range("K32").Value = 3498.121
$ws.Range("M32").Value = -3211.121
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H74").Value = 1710.2307
$ws.Range("I74").Value = 813.3
$ws.Range("K74").Value = 813.3
$ws.Range("M74").Value = 60.70000000000005
$ws.Range("H77").Value = 1710.2307
$ws.Range("I77").Value = 813.3
$ws.Range("K77").Value = 4066.5
$ws.Range("M77").Value = 301.5
$ws.Range("H96").Value = 12316.333
$ws.Range("J96").Value = 12316.333
$ws.Range("L96").Value = 12316.333
$ws.Range("N96").Value = -17808.333
$ws.Range("H100").Value = 34996.668
$ws.Range("J100").Value = 34996.668
$ws.Range("L100").Value = 34996.668
$ws.Range("N100").Value = -37160.668
$ws.Range("H110").Value = 6827.273
$ws.Range("I110").Value = 6888.75
$ws.Range("K110").Value = 6888.75
$ws.Range("M110").Value = -4843.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1697.9166
$ws.Range("I94").Value = 1437.5
$ws.Range("K94").Value = 1437.5
$ws.Range("M94").Value = -986.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 31393.572
$ws.Range("I62").Value = 2902.8572
$ws.Range("K62").Value = 2902.8572
$ws.Range("M62").Value = -2278.8572
$ws.Range("H65").Value = 31393.572
$ws.Range("I65").Value = 2902.8572
$ws.Range("K65").Value = 14514.286
$ws.Range("M65").Value = -11394.286
$ws.Range("H122").Value = 3975.8
$ws.Range("I122").Value = 4737
$ws.Range("K122").Value = 14211
$ws.Range("M122").Value = -11761
$ws.Range("H132").Value = 1263.0605
$ws.Range("I132").Value = 1257.5172
$ws.Range("K132").Value = 3772.5516
$ws.Range("M132").Value = -1242.5516
$ws.Range("H134").Value = 2416.6428
$ws.Range("I134").Value = 2235.3333
$ws.Range("K134").Value = 6705.999899999999
$ws.Range("M134").Value = -4170.999899999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 5000875
$ws.Range("I6").Value = 5000875
$ws.Range("K6").Value = 15002625
$ws.Range("M6").Value = -15002512
$ws.Range("H12").Value = 278.7
$ws.Range("J12").Value = 302.4737
$ws.Range("L12").Value = 907.4211
$ws.Range("N12").Value = -1253.4211
$ws.Range("H97").Value = 10000
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("M98").ClearContents()
$ws.Range("N98").ClearContents()
$ws.Range("H137").Value = 4250.5
$ws.Range("I137").Value = 3495.6667
$ws.Range("K137").Value = 10487.0001
$ws.Range("M137").Value = -5387.000100000001
$ws.Range("H140").Value = 1567.8572
$ws.Range("I140").Value = 995.8333
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 2987.4999
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = 2192.5001
$ws.Range("N140").Value = -25360
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 59999
$ws.Range("J43").Value = 59999
$ws.Range("L43").Value = 59999
$ws.Range("N43").Value = -60301
$ws.Range("H102").Value = 2095.5
$ws.Range("I102").Value = 1934.6
$ws.Range("J102").Value = 2900
$ws.Range("K102").Value = 1934.6
$ws.Range("L102").Value = 2900
$ws.Range("M102").Value = -312.5999999999999
$ws.Range("N102").Value = -6144
$ws.Range("H113").Value = 1627.2
$ws.Range("I113").Value = 1627.2
$ws.Range("K113").Value = 1627.2
$ws.Range("M113").Value = 542.8
$ws.Range("H123").Value = 52000
$ws.Range("J123").Value = 52000
$ws.Range("L123").Value = 52000
$ws.Range("N123").Value = -56900
$ws.Range("H132").Value = 2523.7407
$ws.Range("I132").Value = 2245.6667
$ws.Range("K132").Value = 6737.000100000001
$ws.Range("M132").Value = -4207.000100000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 234
$ws.Range("I16").Value = 250.5
$ws.Range("J16").Value = 201
$ws.Range("K16").Value = 250.5
$ws.Range("L16").Value = 201
$ws.Range("M16").Value = -80.5
$ws.Range("H22").Value = 2926.5
$ws.Range("I22").Value = 1785.0476
$ws.Range("J22").Value = 10916.667
$ws.Range("K22").Value = 1785.0476
$ws.Range("L22").Value = 10916.667
$ws.Range("M22").Value = -1490.0476
$ws.Range("N22").Value = -11506.667
$ws.Range("H27").Value = 2926.5
$ws.Range("I27").Value = 1785.0476
$ws.Range("J27").Value = 10916.667
$ws.Range("K27").Value = 1785.0476
$ws.Range("L27").Value = 10916.667
$ws.Range("M27").Value = -1678.0476
$ws.Range("N27").Value = -11130.667
$ws.Range("H46").Value = 4508.2666
$ws.Range("I46").Value = 2997.25
$ws.Range("J46").Value = 5057.727
$ws.Range("K46").Value = 2997.25
$ws.Range("L46").Value = 5057.727
$ws.Range("M46").Value = -2809.25
$ws.Range("N46").Value = -5433.727
$ws.Range("H55").Value = 829.8889
$ws.Range("I55").Value = 697.8
$ws.Range("J55").Value = 995
$ws.Range("K55").Value = 697.8
$ws.Range("L55").Value = 995
$ws.Range("M55").Value = -524.8
$ws.Range("N55").Value = -1341
$ws.Range("H68").Value = 3545.5557
$ws.Range("I68").Value = 2484
$ws.Range("J68").Value = 5668.6665
$ws.Range("K68").Value = 2484
$ws.Range("L68").Value = 5668.6665
$ws.Range("M68").Value = -1735
$ws.Range("N68").Value = -7166.6665
$ws.Range("H71").Value = 3545.5557
$ws.Range("I71").Value = 2484
$ws.Range("J71").Value = 5668.6665
$ws.Range("K71").Value = 12420
$ws.Range("L71").Value = 28343.3325
$ws.Range("M71").Value = -8676
$ws.Range("N71").Value = -35831.3325
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 2000000
$ws.Range("J26").Value = 2000000
$ws.Range("L26").Value = 2000000
$ws.Range("N26").Value = -2000586
$ws.Range("H96").Value = 2179.6
$ws.Range("I96").Value = 1999.3334
$ws.Range("J96").Value = 2450
$ws.Range("K96").Value = 1999.3334
$ws.Range("L96").Value = 2450
$ws.Range("M96").Value = -626.3334
$ws.Range("N96").Value = -5196
$ws.Range("H107").Value = 1189.1818
$ws.Range("I107").Value = 1091.9375
$ws.Range("K107").Value = 3275.8125
$ws.Range("M107").Value = -1355.8125
$ws.Range("H132").Value = 12259.6
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 14949.5
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 44848.5
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -49908.5
